# Radio module wiring added to the "boat connection" sheet.
# - Adds a "Radio" column (H) note for the D2/D3/D4 rows (7-9), which share
#   the SPI bus/Display connector, and removes the stray "Display/Radio"
#   note that used to sit in C8.
# - Renames the PD2/PD3 interrupt pins to include their INT0/INT1 function,
#   and documents the Radio module's IRQ/CE/CSN pins in columns C and H for
#   rows 22-24.
# - Leaves the final selection on B22, matching where the author ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 7-9: D2/D3/D4 are shared with the Radio module.
$ws.Range("H7").Value = "Radio"

$ws.Range("C8").ClearContents()
$ws.Range("H8").Value = "Radio"

$ws.Range("H9").Value = "Radio"

# Rows 22-24: document the Radio module pin functions.
$ws.Range("A22").Value = "PD2/INT0"
$ws.Range("A23").Value = "PD3/INT1"

$ws.Range("C22").Value = "Radio IRQ"
$ws.Range("H22").Value = "Radio IRQ"

$ws.Range("C23").Value = "Radio CE"
$ws.Range("H23").Value = "Radio CE"

$ws.Range("C24").Value = "Radio CSN"
$ws.Range("H24").Value = "Radio CSN"

# Match the author's final cursor position.
[void]$ws.Range("B22").Select()
